$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "fecha_publicacion" column (E2:E7) currently stores Excel date-serial
# numbers (formatted with a custom "YYYY-MM-DD HH:MM:SS" number format).
# Convert them to literal text values in DD-MM-YYYY form instead.
#
# Pre-formatting the range as Text ("@") stops Excel's automatic
# "this looks like a date" parsing from turning the typed strings back into
# date serials, then resetting the style to Normal afterwards clears the
# number formatting back off of these cells again (matching how the column
# looks once it only holds plain text).
$ws.Range("E2:E7").NumberFormat = "@"

$ws.Range("E2").Value = "10-01-2029"
$ws.Range("E3").Value = "05-04-2030"
$ws.Range("E4").Value = "23-05-2022"
$ws.Range("E5").Value = "10-10-2028"
$ws.Range("E6").Value = "09-09-2024"
$ws.Range("E7").Value = "01-12-2023"

$ws.Range("E2:E7").Style = "Normal"
